$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# --- Update existing rows ---

# Row 2 (TestCase_F1): result SKIP -> PASS
$ws.Range("E2").Value = "PASS"

# Row 4 (TestCase_F3): Jira id + description updated
$ws.Range("B4").Value = "OPQA-210"
$ws.Range("C4").Value = "Verify that user receives a notification when his follower comments on an article contained in his watchlist"

# --- Add new row 7: TestCase_F6 ---
$ws.Range("A7").Value = "TestCase_F6"
$ws.Range("B7").Value = "OPQA-213"
$ws.Range("C7").Value = "Verify that user is able to receive notification when my friend is following some other user."
$ws.Range("D7").Value = "Y"
$ws.Range("E7").Value = "PASS"
$ws.Range("A7:E7").Borders.LineStyle = 1

# --- Add new row 8: TestCase_F7 ---
$ws.Range("A8").Value = "TestCase_F7"
$ws.Range("B8").Value = "OPQA-208"
$ws.Range("C8").Value = "Verify that user receives a notification when someone comments on an article contained in his watchlist"
$ws.Range("D8").Value = "Y"
$ws.Range("E8").Value = "PASS"
$ws.Range("A8:E8").Borders.LineStyle = 1

# --- Update selection to match final workbook state ---
$ws.Range("D6").Select()
